$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new test-scenario row before row 12 (shifts rows 12.. down by one,
# carrying along styles, merged cells and the trailing blank row that becomes
# row 31).
$ws.Rows("12").Insert()

# The insert doesn't reliably keep the A12:F12 cells merged, so redo it.
$null = $ws.Range("A12:F12").Merge()

# Fill in the new question row.
$ws.Range("A12").Value = "Si on rentre une chiffre sup à 9 que cela fonctionne ?"
$ws.Range("G12").Value = "OK"
$ws.Range("H12").Value = "OK"
$ws.Range("I12").Value = "OK"

# Match the saved selection state.
$null = $ws.Range("A13:F13").Select()

# Extend the conditional-formatting ranges to cover the newly added row.
$null = $ws.Range("G2:G30").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G31"))
$null = $ws.Range("H2:H16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H17"))
$null = $ws.Range("I2:I16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I17"))
